# ICAMS-overview.pptx - minor update to powerpoint documentation
#
# 1) Refresh the cached "date last edited" footer field (datetimeFigureOut)
#    from 26/2/2019 to 5/3/2019 everywhere it is defined: once on the
#    slide master and once on each of the 11 slide layouts.
# 2) On slide 2, shrink the red "ReadAndSplitMutectVCFs" label textbox and
#    drop the trailing " (Nanhai to create)" run - the task is done now.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapesOwner, $newDate) {
    $phs = $shapesOwner.Shapes.Placeholders
    for ($i = 1; $i -le $phs.Count; $i++) {
        $ph = $phs.Item($i)
        if ($ph.PlaceholderFormat.Type -eq 16) {
            $ph.TextFrame.TextRange.Text = $newDate
        }
    }
}

$newDate = "5/3/2019"

$master = $p.SlideMaster
Update-DatePlaceholder $master $newDate

$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    Update-DatePlaceholder $layout $newDate
}

# --- Slide 2: "ReadAndSplitMutectVCFs (Nanhai to create)" textbox ---
$slide2 = $p.Slides.Item(2)
$labelShape = $slide2.Shapes.Item("TextBox 39")
$labelShape.TextFrame.TextRange.Text = "ReadAndSplitMutectVCFs"
$labelShape.Width = 207.80473
